$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '23.045.19'
Set-TextValue 'E2' '  -1.02%  '
Set-TextValue 'D3' '1.598.72'
Set-TextValue 'E3' '  -0.06%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  -0.24%  '
Set-TextValue 'E5' '  -0.31%  '
Set-TextValue 'D6' '302.22'
Set-TextValue 'E6' '  +0.01%  '
Set-TextValue 'D7' '0.3778'
Set-TextValue 'E7' '  +0.12%  '
Set-TextValue 'D8' '0.3641'
Set-TextValue 'E8' '  -0.95%  '
Set-TextValue 'D9' '50.76'
Set-TextValue 'E9' '  +1.88%  '
Set-TextValue 'D10' '1.251'
Set-TextValue 'E10' '  -2.23%  '
Set-TextValue 'E11' '  -0.20%  '
Set-TextValue 'D12' '0.08140'
Set-TextValue 'E12' '  +0.14%  '
Set-TextValue 'D13' '22.33'
Set-TextValue 'E13' '  -2.60%  '
Set-TextValue 'D14' '6.571'
Set-TextValue 'E14' '  -1.33%  '
Set-TextValue 'D15' '7.365'
Set-TextValue 'E15' '  -2.89%  '
Set-TextValue 'D16' '0.00001242'
Set-TextValue 'E16' '  -2.14%  '
Set-TextValue 'D17' '1.603.35'
Set-TextValue 'E17' '  +0.33%  '
Set-TextValue 'D18' '91.97'
Set-TextValue 'E18' '  +0.48%  '
Set-TextValue 'D19' '0.06842'
Set-TextValue 'E19' '  +0.25%  '
Set-TextValue 'D20' '18.14'
Set-TextValue 'E20' '  -2.22%  '
Set-TextValue 'D21' '6.509'
Set-TextValue 'E21' '  -1.72%  '
Set-TextValue 'E22' '  -0.23%  '
Set-TextValue 'D23' '13.01'
Set-TextValue 'E23' '  -1.31%  '
Set-TextValue 'D24' '23.050.87'
Set-TextValue 'E24' '  -0.96%  '
Set-TextValue 'D25' '2.363'
Set-TextValue 'E25' '  -0.44%  '
Set-TextValue 'D26' '2.764'
Set-TextValue 'E26' '  -7.54%  '
Set-TextValue 'D27' '21.09'
Set-TextValue 'E27' '  -0.53%  '
Set-TextValue 'D28' '149.16'
Set-TextValue 'E28' '  -1.11%  '
Set-TextValue 'D29' '5.247'
Set-TextValue 'E29' '  -1.48%  '
Set-TextValue 'D30' '134.61'
Set-TextValue 'E30' '  +1.45%  '
Set-TextValue 'D31' '2.356'
Set-TextValue 'E31' '  -4.79%  '
Set-TextValue 'D32' '6.787'
Set-TextValue 'E32' '  -5.25%  '
Set-TextValue 'D33' '1.773.09'
Set-TextValue 'E33' '  -0.08%  '
Set-TextValue 'D34' '0.9582'
Set-TextValue 'E34' '  -1.02%  '
Set-TextValue 'D35' '0.07530'
Set-TextValue 'E35' '  -2.88%  '
Set-TextValue 'B36' 'InternetComputer(DFINITY)'
Set-TextValue 'C36' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D36' '6.198'
Set-TextValue 'E36' '  -1.71%  '
Set-TextValue 'B37' 'VeChain'
Set-TextValue 'C37' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D37' '0.02704'
Set-TextValue 'E37' '  -3.02%  '
Set-TextValue 'B38' 'FraxShare'
Set-TextValue 'C38' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D38' '10.15'
Set-TextValue 'E38' '  -1.10%  '
Set-TextValue 'D39' '0.2518'
Set-TextValue 'E39' '  -1.47%  '
Set-TextValue 'E40' '  -0.74%  '
Set-TextValue 'D41' '1.357'
Set-TextValue 'E41' '  -2.48%  '
Set-TextValue 'D42' '0.7029'
Set-TextValue 'E42' '  -2.18%  '
Set-TextValue 'D43' '12.25'
Set-TextValue 'E43' '  -4.76%  '
Set-TextValue 'D44' '15.19'
Set-TextValue 'E44' '  -6.18%  '
Set-TextValue 'D45' '0.6580'
Set-TextValue 'E45' '  -0.83%  '
Set-TextValue 'B46' 'PancakeSwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D46' '3.998'
Set-TextValue 'E46' '  +0.51%  '
Set-TextValue 'B47' 'NEARProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D47' '2.273'
Set-TextValue 'E47' '  -2.09%  '
Set-TextValue 'B48' 'Quant'
Set-TextValue 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D48' '131.83'
Set-TextValue 'E48' '  -0.17%  '
Set-TextValue 'B49' 'Cronos'
Set-TextValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.07932'
Set-TextValue 'E49' '  -1.04%  '
Set-TextValue 'B50' 'Flow'
Set-TextValue 'C50' 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
Set-TextValue 'D50' '1.221'
Set-TextValue 'E50' '  +3.46%  '
Set-TextValue 'B51' 'ThetaToken'
Set-TextValue 'C51' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D51' '1.226'
Set-TextValue 'E51' '  +2.31%  '

Write-Host "Update complete"
